$d = $word.ActiveDocument

# 1. Trim "the finalized author list" -> "the author list"
$d.Content.Find.Execute(
    "confirm the finalized author list",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "confirm the author list", 2
) | Out-Null

# 2. Rewrite the "Additionally, ..." paragraph: collapse the split runs
#    (which wrapped the now-removed _GoBack bookmark around "0:00") into a
#    single run with the new wording, dropping the specific deadline.
$oldDeadline = "Additionally, should you have any revisions or comments, please communicate them or revise them directly in the overleaf, https://, before Berkeley time, 0:00 11th December. In the absence of further input, I will consider that we have reached consensus, and then I will send you letters for formal approval of the manuscript for submission."
$newDeadline = "Additionally, should you have any revisions or comments, please communicate them or revise them directly in the overleaf, https://, as soon as possible. While other authors might be invited for further revisions later, in the absence of further input, I will consider that we have reached consensus, and then I will send you letters for formal approval of the manuscript for submission."

$d.Content.Find.Execute(
    $oldDeadline,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newDeadline, 2
) | Out-Null

# 3. Re-add the _GoBack bookmark (it was removed along with the runs it used
#    to straddle) onto the last, empty paragraph of the document, which is
#    where Word relocates the auto-maintained _GoBack mark after an edit.
$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastParagraph.Range) | Out-Null
